# Actualización automática: se agrega el nuevo reclamo (fila 76) a la hoja
# "PEBCOM", replicando el mismo formato (texto plano / sin estilo) que el
# resto de las filas de datos ya existentes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 76

# Valores de texto "seguros": Excel no los reinterpreta como numero o fecha,
# asi que se pueden asignar directamente y quedan como texto con el formato
# General por defecto (igual que el resto de las filas).
$ws.Range("C$row").Value = "Brasil 3181"
$ws.Range("F$row").Value = "PEBCOM"
$ws.Range("G$row").Value = "Pendiente"
$ws.Range("H$row").Value = "Picada"
$ws.Range("J$row").Value = "Cambio"
$ws.Range("K$row").Value = "Sin equipos"
$ws.Range("L$row").Value = "Pasante"
$ws.Range("O$row").Value = "San Telmo"
$ws.Range("P$row").Value = "Capital Sur"

# Coordenadas: numericas reales, igual que en el resto de las filas.
$ws.Range("M$row").Value = -58.409002
$ws.Range("N$row").Value = -34.634523

# Valores de texto "ambiguos" (parecen numero o fecha): se fuerza el formato
# de texto antes de escribirlos para que Excel no los convierta, y luego se
# copia el formato "General" de una celda ya existente de esa misma columna
# para que la celda nueva quede sin ningun estilo especial (igual que las
# demas filas de datos, que no tienen atributo de estilo).
$ambiguous = [ordered]@{
    "A" = "-505"
    "B" = "7/11/2025"
    "D" = "4"
    "E" = "808150460"
    "I" = "1"
}

foreach ($col in $ambiguous.Keys) {
    $cellRef = "$col$row"
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $ambiguous[$col]
    $ws.Range("${col}2").Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
}
